# Applies the commit's data edit to TestBattleResults/1_BattlesResults.xlsx
#
# Summary of the change (per the authoritative XML diff):
#  - The "Hornet" and "Zeus" battler names swap places in the shared-string
#    table, so every cell that used to read "Hornet" now reads "Zeus" and
#    vice versa (names only; the row's other stats are addressed below).
#  - A handful of rows have their BT1Type/BT2Type (and the paired
#    Armor/Stamina-before columns) swapped between the "Hornet" and "Zeus"
#    variant of the matchup.
#  - Numerous battle-result metrics (TotalNumberOfRoundsPassed,
#    BT1ArmorAfter, BT2ArmorAfter, BT1StaminaAfter, BT2StaminaAfter,
#    battleWonBT) are updated to new recorded values.
#
# Rather than re-deriving this from game logic, we simply (re)write every
# cell whose final value differs from the "before" workbook, using the
# values dictated by the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = 28
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 681.8
$ws.Range("O2").Value = 2
$ws.Range("J3").Value = 45
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 64.35
$ws.Range("O3").Value = 2
$ws.Range("J4").Value = 27
$ws.Range("K4").Value = 349.83
$ws.Range("M4").Value = 329.3
$ws.Range("C5").Value = "Zeus"
$ws.Range("E5").Value = "GLADIATOR"
$ws.Range("G5").Value = 900
$ws.Range("I5").Value = 499
$ws.Range("J5").Value = 29
$ws.Range("M5").Value = 297.63
$ws.Range("C6").Value = "Hornet"
$ws.Range("E6").Value = "WASP"
$ws.Range("G6").Value = 355
$ws.Range("I6").Value = 599
$ws.Range("J6").Value = 35
$ws.Range("M6").Value = 194.61
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 1
$ws.Range("J7").Value = 18
$ws.Range("J8").Value = 39
$ws.Range("M8").Value = 214.4
$ws.Range("J9").Value = 37
$ws.Range("M9").Value = 544.73
$ws.Range("C10").Value = "Zeus"
$ws.Range("E10").Value = "GLADIATOR"
$ws.Range("G10").Value = 900
$ws.Range("I10").Value = 499
$ws.Range("J10").Value = 21
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 877.71
$ws.Range("N10").Value = 0
$ws.Range("O10").Value = 1
$ws.Range("C11").Value = "Hornet"
$ws.Range("E11").Value = "WASP"
$ws.Range("G11").Value = 355
$ws.Range("I11").Value = 599
$ws.Range("J11").Value = 15
$ws.Range("M11").Value = 234.67
$ws.Range("J12").Value = 18
$ws.Range("M12").Value = 398.66
$ws.Range("J13").Value = 34
$ws.Range("K13").Value = 48.9
$ws.Range("M13").Value = 317.1
$ws.Range("C14").Value = "Zeus"
$ws.Range("E14").Value = "GLADIATOR"
$ws.Range("G14").Value = 900
$ws.Range("I14").Value = 499
$ws.Range("J14").Value = 32
$ws.Range("K14").Value = 15
$ws.Range("M14").Value = 404.49
$ws.Range("C15").Value = "Hornet"
$ws.Range("E15").Value = "WASP"
$ws.Range("G15").Value = 355
$ws.Range("I15").Value = 599
$ws.Range("J15").Value = 28
$ws.Range("M15").Value = 530.46
$ws.Range("J16").Value = 19
$ws.Range("K16").Value = 0
$ws.Range("M16").Value = 172.67
$ws.Range("C17").Value = "Zeus"
$ws.Range("E17").Value = "GLADIATOR"
$ws.Range("G17").Value = 900
$ws.Range("I17").Value = 499
$ws.Range("J17").Value = 38
$ws.Range("N17").Value = 216.92
$ws.Range("C18").Value = "Hornet"
$ws.Range("E18").Value = "WASP"
$ws.Range("G18").Value = 355
$ws.Range("I18").Value = 599
$ws.Range("J18").Value = 32
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 0
$ws.Range("N18").Value = 147.43
$ws.Range("O18").Value = 2
$ws.Range("J19").Value = 35
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = 214.15
$ws.Range("O19").Value = 2
$ws.Range("B20").Value = "Zeus"
$ws.Range("C20").Value = "Hornet"
$ws.Range("D20").Value = "GLADIATOR"
$ws.Range("E20").Value = "WASP"
$ws.Range("F20").Value = 900
$ws.Range("G20").Value = 355
$ws.Range("H20").Value = 499
$ws.Range("I20").Value = 599
$ws.Range("J20").Value = 44
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 24.63
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = 47.53
$ws.Range("O20").Value = 2
$ws.Range("B21").Value = "Zeus"
$ws.Range("D21").Value = "GLADIATOR"
$ws.Range("F21").Value = 900
$ws.Range("H21").Value = 499
$ws.Range("J21").Value = 20
$ws.Range("K21").Value = 138.36
$ws.Range("M21").Value = 499
$ws.Range("N21").Value = 0
$ws.Range("O21").Value = 1
$ws.Range("B22").Value = "Hornet"
$ws.Range("D22").Value = "WASP"
$ws.Range("F22").Value = 900
$ws.Range("H22").Value = 599
$ws.Range("J22").Value = 14
$ws.Range("M22").Value = 308.01

